$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the K-column formula (perp/angular distance log combination):
#    was: 2 * (2 * LOG(E# * SIN(RADIANS(H#)), 2))
#    now: 2 * (LOG(2 * E# * SIN(RADIANS(H#)), 2))
for ($r = 3; $r -le 20; $r++) {
    $ws.Range("K$r").Formula = "=2 * (LOG(2 * E$r * SIN(RADIANS(H$r)), 2))"
}

# 2. Drop the constant factor in M6 from 8 to 0 (lifted restriction mentioned
#    in the commit message).
$ws.Range("M6").Value = 0

# 3. Move/resize the chart (graphicFrame) to its new anchor position.
$co = $ws.ChartObjects().Item(1)
$co.Left = 722.25
$co.Top = 102.75
$co.Width = 433.0625
$co.Height = 216

# 4. Update the active selection on the sheet.
$ws.Activate()
$ws.Range("F6:G27").Select()
